$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet had a redundant "rowZhuyin" column (B) that only ever duplicated
# column A's numeric "row" index, plus a separate "align" column (E) used for
# text alignment of the zhuyin (Bopomofo) column. This edit (adding a PDF
# version of the pinyin/zhuyin guide) cleaned the sheet up:
#   - drop the useless rowZhuyin column
#   - move the "align" column next to "row" (so it reads row, align, chinese,
#     pinyin, zhuyin)
#
# Before: A=row, B=rowZhuyin, C=chinese, D=pinyin, E=align, F=zhuyin
# After:  A=row, B=align,               C=chinese, D=pinyin, E=zhuyin

# 1) Remove the redundant rowZhuyin column.
$ws.Columns.Item(2).Delete() | Out-Null

# Columns are now: A=row, B=chinese, C=pinyin, D=align, E=zhuyin

# 2) Move "align" (now column D) to sit right after "row" (column B),
#    which shifts chinese/pinyin right back into C/D and leaves zhuyin in E.
$ws.Columns.Item(4).Cut() | Out-Null
$ws.Columns.Item(2).Insert() | Out-Null

# The "align" column keeps the width it already had (it just moved from E to
# B, formatting and all). The zhuyin column lands in its new spot (E) with no
# explicit width override, same as it had none in F before — so give it a
# fresh best-fit width now that it's the sheet's last column.
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

$ws.Range("I18").Select() | Out-Null
